{"js": "// Insert a new paragraph right after the paragraph that tells the model to\n// write \"\u0441\u0442\u0430\u0442\u0443\u0441 \u043e\u0436\u0438\u0434\u0430\u0435\u0442 \u0437\u0432\u043e\u043d\u043a\u0430\", containing the new \"\u043d\u0435\u0443\u0441\u043f\u0435\u0448\u043d\u044b\u0439 \u0434\u0438\u0430\u043b\u043e\u0433\"\n// instruction. The new paragraph mirrors the formatting of the paragraph it\n// follows (same fonts / size / spacing), matching the author's edit.\n\nconst anchorText = \"\u043f\u043e\u0441\u043b\u0435\u0434\u043d\u0438\u043c \u043f\u0440\u0435\u0434\u043b\u043e\u0436\u0435\u043d\u0438\u0435\u043c \u043d\u0430\u043f\u0438\u0448\u0438: \\\"\u0441\u0442\u0430\u0442\u0443\u0441 \u043e\u0436\u0438\u0434\u0430\u0435\u0442 \u0437\u0432\u043e\u043d\u043a\u0430\\\"\";\nconst newText =\n  \"\u0410\u043d\u0430\u043b\u0438\u0437\u0438\u0440\u0443\u044f \u0432\u0435\u0441\u044c \u043a\u043e\u043d\u0442\u0435\u043a\u0441\u0442 \u0434\u0438\u0430\u043b\u043e\u0433\u0430, \u0435\u0441\u043b\u0438 \u0442\u044b \u0443\u0431\u0435\u0436\u0434\u0430\u0435\u0448\u044c\u0441\u044f, \u0447\u0442\u043e \u043a\u043b\u0438\u0435\u043d\u0442 \u0434\u0430\u043b \" +\n  \"\u043e\u0442\u0440\u0438\u0446\u0430\u0442\u0435\u043b\u044c\u043d\u044b\u0439 \u043e\u0442\u0432\u0435\u0442 \u043d\u0430 \u043f\u0440\u0435\u0434\u043b\u043e\u0436\u0435\u043d\u0438\u0435 \u0441\u043e\u0437\u0432\u043e\u043d\u0438\u0442\u044c\u0441\u044f \u0438\u043b\u0438 \u0441\u043a\u0430\u0437\u0430\u043b, \u0447\u0442\u043e \u043d\u0438\u0447\u0435\u0433\u043e \" +\n  \"\u043d\u0435 \u043d\u0430\u0434\u043e, \u0442\u043e \u043f\u043e\u0441\u043b\u0435\u0434\u043d\u0438\u043c \u043f\u0440\u0435\u0434\u043b\u043e\u0436\u0435\u043d\u0438\u0435\u043c \u043d\u0430\u043f\u0438\u0448\u0438: \\\"\u043d\u0435\u0443\u0441\u043f\u0435\u0448\u043d\u044b\u0439 \u0434\u0438\u0430\u043b\u043e\u0433\\\"\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet anchorParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text && p.text.indexOf(anchorText) !== -1) {\n    anchorParagraph = p;\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error(\"Could not find the anchor paragraph ('\u0441\u0442\u0430\u0442\u0443\u0441 \u043e\u0436\u0438\u0434\u0430\u0435\u0442 \u0437\u0432\u043e\u043d\u043a\u0430').\");\n}\n\n// Insert the new paragraph directly after the anchor; Word clones the\n// anchor paragraph's mark formatting (fonts, size, spacing, alignment) for\n// the new paragraph automatically, matching the sibling paragraphs already\n// in this prompt document.\nanchorParagraph.insertParagraph(newText, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert a new paragraph right after the paragraph that instructs the model\n# to write \"\u0441\u0442\u0430\u0442\u0443\u0441 \u043e\u0436\u0438\u0434\u0430\u0435\u0442 \u0437\u0432\u043e\u043d\u043a\u0430\", adding the new \"\u043d\u0435\u0443\u0441\u043f\u0435\u0448\u043d\u044b\u0439 \u0434\u0438\u0430\u043b\u043e\u0433\"\n# instruction. The new paragraph inherits the formatting (fonts, size,\n# spacing, alignment) of the paragraph it follows, matching the sibling\n# paragraphs already in this prompt document.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"\u043f\u043e\u0441\u043b\u0435\u0434\u043d\u0438\u043c \u043f\u0440\u0435\u0434\u043b\u043e\u0436\u0435\u043d\u0438\u0435\u043c \u043d\u0430\u043f\u0438\u0448\u0438: \" + [char]34 + \"\u0441\u0442\u0430\u0442\u0443\u0441 \u043e\u0436\u0438\u0434\u0430\u0435\u0442 \u0437\u0432\u043e\u043d\u043a\u0430\" + [char]34\n$newText = \"\u0410\u043d\u0430\u043b\u0438\u0437\u0438\u0440\u0443\u044f \u0432\u0435\u0441\u044c \u043a\u043e\u043d\u0442\u0435\u043a\u0441\u0442 \u0434\u0438\u0430\u043b\u043e\u0433\u0430, \u0435\u0441\u043b\u0438 \u0442\u044b \u0443\u0431\u0435\u0436\u0434\u0430\u0435\u0448\u044c\u0441\u044f, \u0447\u0442\u043e \u043a\u043b\u0438\u0435\u043d\u0442 \u0434\u0430\u043b \u043e\u0442\u0440\u0438\u0446\u0430\u0442\u0435\u043b\u044c\u043d\u044b\u0439 \u043e\u0442\u0432\u0435\u0442 \u043d\u0430 \u043f\u0440\u0435\u0434\u043b\u043e\u0436\u0435\u043d\u0438\u0435 \u0441\u043e\u0437\u0432\u043e\u043d\u0438\u0442\u044c\u0441\u044f \u0438\u043b\u0438 \u0441\u043a\u0430\u0437\u0430\u043b, \u0447\u0442\u043e \u043d\u0438\u0447\u0435\u0433\u043e \u043d\u0435 \u043d\u0430\u0434\u043e, \u0442\u043e \u043f\u043e\u0441\u043b\u0435\u0434\u043d\u0438\u043c \u043f\u0440\u0435\u0434\u043b\u043e\u0436\u0435\u043d\u0438\u0435\u043c \u043d\u0430\u043f\u0438\u0448\u0438: \" + [char]34 + \"\u043d\u0435\u0443\u0441\u043f\u0435\u0448\u043d\u044b\u0439 \u0434\u0438\u0430\u043b\u043e\u0433\" + [char]34\n\n$anchorParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($anchorText)) {\n        $anchorParagraph = $p\n        break\n    }\n}\n\nif ($anchorParagraph -eq $null) {\n    throw \"Could not find the anchor paragraph ('\u0441\u0442\u0430\u0442\u0443\u0441 \u043e\u0436\u0438\u0434\u0430\u0435\u0442 \u0437\u0432\u043e\u043d\u043a\u0430').\"\n}\n\n$anchorParagraph.Range.InsertParagraphAfter()\n$newParagraph = $anchorParagraph.Next()\n$newParagraph.Range.Text = $newText\n"}
